$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N1").Value = "Time"
$ws.Range("O1").Value = "Owner"

$ws.Range("O1").Select()
